# Scheduled market-price refresh: Universalis price pull updates the
# computed Leve profit columns (H:N) for the rows whose items repriced.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 230.94444
$ws.Range("I4").Value = 236.41667
$ws.Range("K4").Value = 236.41667
$ws.Range("M4").Value = -122.41667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2333
$ws.Range("I19").Value = 5441.6
$ws.Range("J19").Value = 696.8946999999999
$ws.Range("K19").Value = 5441.6
$ws.Range("L19").Value = 696.8946999999999
$ws.Range("M19").Value = -5266.6
$ws.Range("N19").Value = -1046.8947

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1831.8125
$ws.Range("J43").Value = 1693
$ws.Range("L43").Value = 1693
$ws.Range("N43").Value = -1831

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 267
$ws.Range("I61").Value = 111.666664
$ws.Range("J61").Value = 500
$ws.Range("K61").Value = 334.999992
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -162.999992
$ws.Range("N61").Value = -1844

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2814.6428
$ws.Range("I116").Value = 3101.625
$ws.Range("J116").Value = 2432
$ws.Range("K116").Value = 3101.625
$ws.Range("L116").Value = 2432
$ws.Range("M116").Value = 340.375
$ws.Range("N116").Value = -9316

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2251
$ws.Range("I132").Value = 1466.26
$ws.Range("J132").Value = 5269.231
$ws.Range("K132").Value = 4398.78
$ws.Range("L132").Value = 15807.693
$ws.Range("M132").Value = -1868.78
$ws.Range("N132").Value = -20867.693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1812.7646
$ws.Range("I2").Value = 1812.7646
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1812.7646
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1699.7646
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5404.769
$ws.Range("I63").Value = 5404.769
$ws.Range("K63").Value = 5404.769
$ws.Range("M63").Value = -4718.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 5404.769
$ws.Range("I66").Value = 5404.769
$ws.Range("K66").Value = 27023.845
$ws.Range("M66").Value = -23591.845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1052.9474
$ws.Range("I97").Value = 1087.9412
$ws.Range("J97").Value = 755.5
$ws.Range("K97").Value = 1087.9412
$ws.Range("L97").Value = 755.5
$ws.Range("M97").Value = -591.9412
$ws.Range("N97").Value = -1747.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1251.1708
$ws.Range("I110").Value = 1156.2188
$ws.Range("J110").Value = 1588.7778
$ws.Range("K110").Value = 1156.2188
$ws.Range("L110").Value = 1588.7778
$ws.Range("M110").Value = 888.7811999999999
$ws.Range("N110").Value = -5678.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1812.7646
$ws.Range("I116").Value = 1812.7646
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1812.7646
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 481.2354
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 21661.246
$ws.Range("I132").Value = 27518.7
$ws.Range("J132").Value = 3638.3076
$ws.Range("K132").Value = 82556.10000000001
$ws.Range("L132").Value = 10914.9228
$ws.Range("M132").Value = -80026.10000000001
$ws.Range("N132").Value = -15974.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1812.7646
$ws.Range("I3").Value = 1812.7646
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1812.7646
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1698.7646
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 950.6774
$ws.Range("I20").Value = 764.1111
$ws.Range("J20").Value = 1209
$ws.Range("K20").Value = 764.1111
$ws.Range("L20").Value = 1209
$ws.Range("M20").Value = -517.1111
$ws.Range("N20").Value = -1703

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7475
$ws.Range("I99").Value = 13250
$ws.Range("J99").Value = 1700
$ws.Range("K99").Value = 13250
$ws.Range("L99").Value = 1700
$ws.Range("M99").Value = -11752
$ws.Range("N99").Value = -4696

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 964.91113
$ws.Range("I107").Value = 923.9143
$ws.Range("J107").Value = 1108.4
$ws.Range("K107").Value = 923.9143
$ws.Range("L107").Value = 1108.4
$ws.Range("M107").Value = 996.0857
$ws.Range("N107").Value = -4948.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 73355.14
$ws.Range("I99").Value = 112411.555
$ws.Range("J99").Value = 3053.6
$ws.Range("K99").Value = 112411.555
$ws.Range("L99").Value = 3053.6
$ws.Range("M99").Value = -110913.555
$ws.Range("N99").Value = -6049.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 816.2963
$ws.Range("I105").Value = 746.9048
$ws.Range("K105").Value = 746.9048
$ws.Range("M105").Value = 1000.0952

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 73355.14
$ws.Range("I126").Value = 112411.555
$ws.Range("J126").Value = 3053.6
$ws.Range("K126").Value = 337234.665
$ws.Range("L126").Value = 9160.799999999999
$ws.Range("M126").Value = -334764.665
$ws.Range("N126").Value = -14100.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 28299.666
$ws.Range("J141").Value = 28299.666
$ws.Range("L141").Value = 28299.666
$ws.Range("N141").Value = -38659.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1083.575
$ws.Range("I5").Value = 429.90475
$ws.Range("K5").Value = 1289.71425
$ws.Range("M5").Value = -1177.71425

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 39.21875
$ws.Range("I12").Value = 24.25
$ws.Range("J12").Value = 44.208332
$ws.Range("K12").Value = 72.75
$ws.Range("L12").Value = 132.624996
$ws.Range("M12").Value = 100.25
$ws.Range("N12").Value = -478.624996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3052.9412
$ws.Range("I39").Value = 666.6667
$ws.Range("J39").Value = 3564.2856
$ws.Range("K39").Value = 2000.0001
$ws.Range("L39").Value = 10692.8568
$ws.Range("M39").Value = -1706.0001
$ws.Range("N39").Value = -11280.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 572.55554
$ws.Range("I107").Value = 198.27272
$ws.Range("J107").Value = 829.875
$ws.Range("K107").Value = 594.81816
$ws.Range("L107").Value = 2489.625
$ws.Range("M107").Value = 1325.18184
$ws.Range("N107").Value = -6329.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 5087.7334
$ws.Range("I116").Value = 310
$ws.Range("K116").Value = 930
$ws.Range("M116").Value = 2512

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 589.2
$ws.Range("I122").Value = 432.27777
$ws.Range("J122").Value = 992.7143
$ws.Range("K122").Value = 3890.49993
$ws.Range("L122").Value = 8934.4287
$ws.Range("M122").Value = -1440.49993
$ws.Range("N122").Value = -13834.4287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1083.575
$ws.Range("I135").Value = 429.90475
$ws.Range("K135").Value = 3869.14275
$ws.Range("M135").Value = -1334.14275

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 715
$ws.Range("I97").Value = 654.2857
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 654.2857
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -158.2857
$ws.Range("N97").Value = -1792

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2171.6296
$ws.Range("I126").Value = 1937.4546
$ws.Range("J126").Value = 3202
$ws.Range("K126").Value = 5812.3638
$ws.Range("L126").Value = 9606
$ws.Range("M126").Value = -3342.3638
$ws.Range("N126").Value = -14546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3555.923
$ws.Range("I132").Value = 3391.6875
$ws.Range("J132").Value = 3818.7
$ws.Range("K132").Value = 10175.0625
$ws.Range("L132").Value = 11456.1
$ws.Range("M132").Value = -7645.0625
$ws.Range("N132").Value = -16516.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 587.53845
$ws.Range("I22").Value = 417.91666
$ws.Range("J22").Value = 732.9286
$ws.Range("K22").Value = 417.91666
$ws.Range("L22").Value = 732.9286
$ws.Range("M22").Value = -122.91666
$ws.Range("N22").Value = -1322.9286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 587.53845
$ws.Range("I27").Value = 417.91666
$ws.Range("J27").Value = 732.9286
$ws.Range("K27").Value = 417.91666
$ws.Range("L27").Value = 732.9286
$ws.Range("M27").Value = -310.91666
$ws.Range("N27").Value = -946.9286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1621.5834
$ws.Range("I93").Value = 1651
$ws.Range("J93").Value = 1533.3334
$ws.Range("K93").Value = 1651
$ws.Range("L93").Value = 1533.3334
$ws.Range("M93").Value = -403
$ws.Range("N93").Value = -4029.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 83339240
$ws.Range("I100").Value = 8528.571
$ws.Range("K100").Value = 8528.571
$ws.Range("M100").Value = -7987.571
